$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New opportunity data for rows 2-16 (columns A-H)
$data = @(
    @("1331862", "https://aiesec.org/opportunity/global-talent/1331862", "Communication Analyst", "Buenos Aires, Cdad. Autónoma de Buenos Aires, Argentina", "No", "0 applicants", "3 - 6 Months", "Beckett Idiomas SA"),
    @("1331814", "https://aiesec.org/opportunity/global-talent/1331814", "Digital Marketing Coordinator", "Unawatuna, Sri Lanka", "No", "0 applicants", "3 - 6 Months", "Angel Beach Club Unawatuna"),
    @("1331811", "https://aiesec.org/opportunity/global-talent/1331811", "[EXP] Global Digital Marketing Intern – Content & Social Media", "Maastricht, Netherlands", "Yes", "15 applicants", "6 - 18 Months", "DHL Group"),
    @("1331777", "https://aiesec.org/opportunity/global-talent/1331777", "Social Media Manager", "Hong Kong", "No", "11 applicants", "6 - 18 Months", "Wong's Limited"),
    @("1331737", "https://aiesec.org/opportunity/global-talent/1331737", "Accelerate Romania - Electronics Engineer. Industrial Equipment Repair", "Cluj-Napoca, Romania", "No", "1 applicant", "9 - 12 Weeks", "Apelogic Engineering"),
    @("1331692", "https://aiesec.org/opportunity/global-talent/1331692", "Accelerate Romania - Workshop Supervisor & Industrial Equipment Repair", "Cluj-Napoca, Romania", "No", "1 applicant", "6 - 18 Months", "Apelogic Engineering"),
    @("1331690", "https://aiesec.org/opportunity/global-talent/1331690", "Accelerate Romania - Account Manager", "Cluj-Napoca, Romania", "No", "1 applicant", "6 - 18 Months", "Apelogic Engineering"),
    @("1331662", "https://aiesec.org/opportunity/global-talent/1331662", "Mechanical Engineering Intern", "Phagwara, Punjab, India", "No", "0 applicants", "3 - 6 Months", "GNA University"),
    @("1331430", "https://aiesec.org/opportunity/global-talent/1331430", "Digital Media Strategist - Mid Term", "Nugegoda, Sri Lanka", "No", "0 applicants", "3 - 6 Months", "Brand Corridor (Pvt) Ltd"),
    @("1331393", "https://aiesec.org/opportunity/global-talent/1331393", "Back end developer", "Jawhara, Tunisie", "No", "4 applicants", "9 - 12 Weeks", "Digital Glow Agency"),
    @("1327293", "https://aiesec.org/opportunity/global-talent/1327293", "Business Executive", "Hong Kong", "No", "108 applicants", "6 - 18 Months", "ASA Building Materials (HK) Limited"),
    @("1322725", "https://aiesec.org/opportunity/global-talent/1322725", "Software Developer", "33 Bielefeld, Deutschland", "No", "102 applicants", "6 - 18 Months", "Software Development  Project Internship"),
    @("1315190", "https://aiesec.org/opportunity/global-talent/1315190", "Front Office Trainee", "Hong Kong", "No", "62 applicants", "6 - 18 Months", "Park Hotel International Limited"),
    @("1315102", "https://aiesec.org/opportunity/global-talent/1315102", "Food and Beverage Trainee", "Hong Kong", "No", "25 applicants", "6 - 18 Months", "Park Hotel International Limited"),
    @("1301518", "https://aiesec.org/opportunity/global-talent/1301518", "MARKETING", "Yıldırım, Türkiye", "No", "85 applicants", "9 - 12 Weeks", "OMTEC Automotive")
)

# The "PREMIUM = Yes" highlight moves from row 3 to row 4 - copy the cell
# formats (not values) so the existing highlight style is reused verbatim
# rather than minting new style/fill entries.
$ws.Range("E3").Copy()
$ws.Range("E4").PasteSpecial(-4122)
$ws.Range("E2").Copy()
$ws.Range("E3").PasteSpecial(-4122)

# Remove the old rows (17-23) that no longer exist in the refreshed scrape
$ws.Range("A17:H23").EntireRow.Delete()

# Opportunity IDs (col A) are textual identifiers, not numbers - keep them as text
$ws.Range("A2:A16").NumberFormat = "@"

for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $i + 2
    $row = $data[$i]
    for ($j = 0; $j -lt $row.Length; $j++) {
        $c = $j + 1
        $ws.Cells.Item($r, $c).Value = $row[$j]
    }
}

# Column width adjustments
$ws.Range("C1").EntireColumn.ColumnWidth = 73
$ws.Range("D1").EntireColumn.ColumnWidth = 58
$ws.Range("H1").EntireColumn.ColumnWidth = 43
